# 생치 column chart 추가
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new hospital rows (18-20)
$ws.Range("A18").Value = "성빈센트"
$ws.Range("B18").Value = 37.278156000000003
$ws.Range("C18").Value = 127.02781

$ws.Range("A19").Value = "현대병원"
$ws.Range("B19").Value = 37.715788000000003
$ws.Range("C19").Value = 127.179631

$ws.Range("A20").Value = "세종병원"
$ws.Range("B20").Value = 37.481034000000001
$ws.Range("C20").Value = 126.791188

# Update selection to match final state
$ws.Range("D8").Select()
